$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 582.5
$ws.Range("I33").Value = 228.125
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 228.125
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = 0.875
$ws.Range("N33").Value = -2458

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4679.5386
$ws.Range("I40").Value = 3555.647
$ws.Range("J40").Value = 6802.4443
$ws.Range("K40").Value = 3555.647
$ws.Range("L40").Value = 6802.4443
$ws.Range("M40").Value = -3380.647
$ws.Range("N40").Value = -7152.4443

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1696.8438
$ws.Range("I32").Value = 599.6957
$ws.Range("J32").Value = 4500.6665
$ws.Range("K32").Value = 599.6957
$ws.Range("L32").Value = 4500.6665
$ws.Range("M32").Value = -312.6957
$ws.Range("N32").Value = -5074.6665

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4037.625
$ws.Range("I102").Value = 2050.1667
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 2050.1667
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -428.1667000000002
$ws.Range("N102").Value = -13244

# ARM row 106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 20000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 20000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 20000
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -22524

# ARM row 119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 34500
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 34500
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 34500
$ws.Range("N119").Value = -44176

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4281.4736
$ws.Range("I132").Value = 3646.8125
$ws.Range("J132").Value = 7666.3335
$ws.Range("K132").Value = 10940.4375
$ws.Range("L132").Value = 22999.0005
$ws.Range("M132").Value = -8410.4375
$ws.Range("N132").Value = -28059.0005

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 32499.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 32499.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 32499.5
$ws.Range("N135").Value = -42639.5

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1999.4
$ws.Range("I20").Value = 1999.4
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1999.4
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1752.4

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2037.7778
$ws.Range("I99").Value = 2091.6667
$ws.Range("J99").Value = 1930
$ws.Range("K99").Value = 2091.6667
$ws.Range("L99").Value = 1930
$ws.Range("M99").Value = -593.6667000000002
$ws.Range("N99").Value = -4926

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3160.9
$ws.Range("I105").Value = 3229.375
$ws.Range("J105").Value = 2887
$ws.Range("K105").Value = 3229.375
$ws.Range("L105").Value = 2887
$ws.Range("M105").Value = -1482.375
$ws.Range("N105").Value = -6381

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2444
$ws.Range("I58").Value = 1900.875
$ws.Range("J58").Value = 3168.1667
$ws.Range("K58").Value = 1900.875
$ws.Range("L58").Value = 3168.1667
$ws.Range("M58").Value = -1697.875
$ws.Range("N58").Value = -3574.1667

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2078.1177
$ws.Range("I105").Value = 1960.6666
$ws.Range("J105").Value = 2360
$ws.Range("K105").Value = 1960.6666
$ws.Range("L105").Value = 2360
$ws.Range("M105").Value = -213.6666
$ws.Range("N105").Value = -5854

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 429.25
$ws.Range("I107").Value = 294.4
$ws.Range("J107").Value = 525.5714
$ws.Range("K107").Value = 294.4
$ws.Range("L107").Value = 525.5714
$ws.Range("M107").Value = 1625.6
$ws.Range("N107").Value = -4365.5714

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2444
$ws.Range("I136").Value = 1900.875
$ws.Range("J136").Value = 3168.1667
$ws.Range("K136").Value = 5702.625
$ws.Range("L136").Value = 9504.500100000001
$ws.Range("M136").Value = -3152.625
$ws.Range("N136").Value = -14604.5001

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 802.5
$ws.Range("I5").Value = 746
$ws.Range("J5").Value = 859
$ws.Range("K5").Value = 2238
$ws.Range("L5").Value = 2577
$ws.Range("M5").Value = -2126
$ws.Range("N5").Value = -2801

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 802.5
$ws.Range("I135").Value = 746
$ws.Range("J135").Value = 859
$ws.Range("K135").Value = 6714
$ws.Range("L135").Value = 7731
$ws.Range("M135").Value = -4179
$ws.Range("N135").Value = -12801

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 342.75
$ws.Range("I2").Value = 106
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 106
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = -2226

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6294.067
$ws.Range("I7").Value = 5619.143
$ws.Range("J7").Value = 6884.625
$ws.Range("K7").Value = 5619.143
$ws.Range("L7").Value = 6884.625
$ws.Range("M7").Value = -5507.143
$ws.Range("N7").Value = -7108.625

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2017
$ws.Range("I22").Value = 1525.5
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1525.5
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1230.5
$ws.Range("N22").Value = -3590

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2017
$ws.Range("I27").Value = 1525.5
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 1525.5
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -1418.5
$ws.Range("N27").Value = -3214

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2696.5789
$ws.Range("I61").Value = 1282.4
$ws.Range("J61").Value = 7999.75
$ws.Range("K61").Value = 1282.4
$ws.Range("L61").Value = 7999.75
$ws.Range("M61").Value = -1080.4
$ws.Range("N61").Value = -8403.75

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7562
$ws.Range("I68").Value = 3248
$ws.Range("J68").Value = 9000
$ws.Range("K68").Value = 3248
$ws.Range("L68").Value = 9000
$ws.Range("M68").Value = -2499
$ws.Range("N68").Value = -10498

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 7562
$ws.Range("I71").Value = 3248
$ws.Range("J71").Value = 9000
$ws.Range("K71").Value = 16240
$ws.Range("L71").Value = 45000
$ws.Range("M71").Value = -12496
$ws.Range("N71").Value = -52488

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6376.769
$ws.Range("I100").Value = 4065.6667
$ws.Range("J100").Value = 8357.714
$ws.Range("K100").Value = 4065.6667
$ws.Range("L100").Value = 8357.714
$ws.Range("M100").Value = -3524.6667
$ws.Range("N100").Value = -9439.714

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2696.5789
$ws.Range("I113").Value = 1282.4
$ws.Range("J113").Value = 7999.75
$ws.Range("K113").Value = 1282.4
$ws.Range("L113").Value = 7999.75
$ws.Range("M113").Value = 887.5999999999999
$ws.Range("N113").Value = -12339.75

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6294.067
$ws.Range("I126").Value = 5619.143
$ws.Range("J126").Value = 6884.625
$ws.Range("K126").Value = 16857.429
$ws.Range("L126").Value = 20653.875
$ws.Range("M126").Value = -14387.429
$ws.Range("N126").Value = -25593.875

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3667.2727
$ws.Range("I136").Value = 2620
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 7860
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -5310
$ws.Range("N136").Value = -21600

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11750.5
$ws.Range("I62").Value = 9001
$ws.Range("J62").Value = 14500
$ws.Range("K62").Value = 9001
$ws.Range("L62").Value = 14500
$ws.Range("M62").Value = -8377
$ws.Range("N62").Value = -15748

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 11750.5
$ws.Range("I65").Value = 9001
$ws.Range("J65").Value = 14500
$ws.Range("K65").Value = 45005
$ws.Range("L65").Value = 72500
$ws.Range("M65").Value = -41885
$ws.Range("N65").Value = -78740

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2667.72
$ws.Range("I136").Value = 1983.8422
$ws.Range("J136").Value = 4833.3335
$ws.Range("K136").Value = 5951.5266
$ws.Range("L136").Value = 14500.0005
$ws.Range("M136").Value = -3401.5266
$ws.Range("N136").Value = -19600.0005
